$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the "datetimeFigureOut" date placeholder field text on the slide
#    master and every slide layout: 12/1/2012 -> 12/8/2012
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "12/1/2012") {
                $shp.TextFrame.TextRange.Text = "12/8/2012"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $lay = $master.CustomLayouts.Item($i)
    Update-DateShapes $lay.Shapes
}

# ---------------------------------------------------------------------------
# 2. Diagram edits on slide 1
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# Reposition "Rectangle 181" (TeamData) and "Rectangle 182" (EvalResultData)
# boxes down slightly (target offsets are 3322817 EMU and 3751442 EMU; a tiny
# epsilon is added in points so the float->EMU conversion lands exactly on
# the target instead of truncating one EMU short).
$teamData = $s.Shapes.Item("Rectangle 181")
$teamData.Top = 261.6391738582677

$evalResultData = $s.Shapes.Item("Rectangle 182")
$evalResultData.Top = 295.3891738582678

# Toggling the shadow off stamps an explicit (empty) effect list onto these
# connectors/boxes, matching the diagram cleanup.
$effectShapeNames = @(
    "Elbow Connector 167",
    "Rectangle 171",
    "Rectangle 178"
)
foreach ($name in $effectShapeNames) {
    $shp = $s.Shapes.Item($name)
    $shp.Shadow.Visible = $false
}

# The remaining four "Elbow Connector 68" shapes share a name, so address
# them positionally (z-order) instead.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Elbow Connector 68") {
        $shp.Shadow.Visible = $false
    }
}

# Delete the now-unused "Rounded Rectangle 156" background shape.
$s.Shapes.Item("Rounded Rectangle 156").Delete()
